# Refresh the crypto price/volume table (GitHub Actions data refresh).
# For D-column cells whose new value looks like a plain number (e.g. "214.11"),
# prefix with an apostrophe so Excel stores it as text (matching the original
# inlineStr cells) instead of silently converting it to a numeric value; then
# reset the cell Style back to "Normal" so no stray number-format style lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.050.97"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "1.649.46"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'214.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'23.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.40%  "
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").Value = "1.883.45"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").Value = "1.648.25"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").Value = "'65.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "28.053.67"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("D18").Value = "'232.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'10.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.81%  "
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D25").Value = "'152.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").Value = "'15.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").Value = "1.447.44"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'0.894"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").Value = "'0.931"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.83%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "'69.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("E45").Value = "  +5.67%  "
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("D48").Value = "1.791.92"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("D49").Value = "'89.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  +0.21%  "
